# Apply the NATMI recompute edit (Thbs1-Itga4) per commit "Natmi following Dr Hou advice"
# Recomputes ligand/receptor expression stats now that the sending/target
# clusters have 2-3 expressing cells each (was 1), expanding rows 2-31 into
# rows 2-37 (adds the "sCs" sending-cluster block as rows 32-37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Thbs1", "Itga4", "ECs", 2, 1, 31.645482, 63.290964, 0.03555980726701226, 0.0244846141215985, 2, 1, 0.317263, 0.634526, 0.002587741009497455, 0.001733186563985812, 10.039940555766, 40.15976222306401, 0.00009201957155467325, 0.0000424364042199318),
    @("ECs", "Thbs1", "Itga4", "FAPs", 2, 1, 31.645482, 63.290964, 0.03555980726701226, 0.0244846141215985, 1, 0.3333333333333333, 0.02842633333333333, 0.085279, 0.0002318580752130288, 0.0002329367385893503, 0.899565019826, 5.397390118955999, 0.000008244828467875736, 0.000005703366159103904),
    @("ECs", "Thbs1", "Itga4", "M1", 2, 1, 31.645482, 63.290964, 0.03555980726701226, 0.0244846141215985, 3, 1, 49.39115366666666, 148.173461, 0.4028566641859401, 0.4047308569616938, 1563.006864317734, 9378.041185906402, 0.01432550533468351, 0.00990967885581095),
    @("ECs", "Thbs1", "Itga4", "M2", 2, 1, 31.645482, 63.290964, 0.03555980726701226, 0.0244846141215985, 3, 1, 70.865733, 212.597199, 0.5780130789036135, 0.5807021443531359, 2242.580278068306, 13455.48166840984, 0.02055403368362485, 0.01421826792407132),
    @("ECs", "Thbs1", "Itga4", "Neutro", 2, 1, 31.645482, 63.290964, 0.03555980726701226, 0.0244846141215985, 3, 1, 0.6137786666666667, 1.841336, 0.005006257352694774, 0.005029547749002215, 19.423321747984, 116.539930487904, 0.0001780215465908892, 0.0001231465358404736),
    @("ECs", "Thbs1", "Itga4", "sCs", 2, 1, 31.645482, 63.290964, 0.03555980726701226, 0.0244846141215985, 2, 1, 1.3859455, 2.771891, 0.01130440047304115, 0.007571327633592945, 43.858913373231, 175.435653492924, 0.0004019823020904656, 0.0001853810354967188),
    @("FAPs", "Thbs1", "Itga4", "ECs", 3, 1, 83.22744366666667, 249.682331, 0.09352209759714789, 0.09659160077758068, 2, 1, 0.317263, 0.634526, 0.002587741009497455, 0.001733186563985812, 26.40498846001767, 158.429930760106, 0.000242010967246363, 0.0001674112646615844),
    @("FAPs", "Thbs1", "Itga4", "FAPs", 3, 1, 83.22744366666667, 249.682331, 0.09352209759714789, 0.09659160077758068, 1, 0.3333333333333333, 0.02842633333333333, 0.085279, 0.0002318580752130288, 0.0002329367385893503, 2.365851056149889, 21.292659505349, 0.00002168385353875973, 0.00002249973246025419),
    @("FAPs", "Thbs1", "Itga4", "M1", 3, 1, 83.22744366666667, 249.682331, 0.09352209759714789, 0.09659160077758068, 3, 1, 49.39115366666666, 148.173461, 0.4028566641859401, 0.4047308569616938, 4110.699459424177, 36996.29513481758, 0.03767600026565892, 0.03909360135801204),
    @("FAPs", "Thbs1", "Itga4", "M2", 3, 1, 83.22744366666667, 249.682331, 0.09352209759714789, 0.09659160077758068, 3, 1, 70.865733, 212.597199, 0.5780130789036135, 0.5807021443531359, 5897.973801154542, 53081.76421039087, 0.05405699557765169, 0.05609094969804313),
    @("FAPs", "Thbs1", "Itga4", "Neutro", 3, 1, 83.22744366666667, 249.682331, 0.09352209759714789, 0.09659160077758068, 3, 1, 0.6137786666666667, 1.841336, 0.005006257352694774, 0.005029547749002215, 51.08322940380178, 459.749064634216, 0.0004681956887351598, 0.0004858120682634015),
    @("FAPs", "Thbs1", "Itga4", "sCs", 3, 1, 83.22744366666667, 249.682331, 0.09352209759714789, 0.09659160077758068, 2, 1, 1.3859455, 2.771891, 0.01130440047304115, 0.007571327633592945, 115.3487010263202, 692.092206157921, 0.001057211244316999, 0.0007313266561402744),
    @("M1", "Thbs1", "Itga4", "ECs", 3, 1, 524.5768889999999, 1573.730667, 0.5894633891046084, 0.6088102578564109, 2, 1, 0.317263, 0.634526, 0.002587741009497455, 0.001733186563985812, 166.428837534807, 998.573025208842, 0.00152537858558335, 0.001055181758933469),
    @("M1", "Thbs1", "Itga4", "FAPs", 3, 1, 524.5768889999999, 1573.730667, 0.5894633891046084, 0.6088102578564109, 1, 0.3333333333333333, 0.02842633333333333, 0.085279, 0.0002318580752130288, 0.0002329367385893503, 14.911797505677, 134.206177551093, 0.0001366718468063431, 0.0001418142758848137),
    @("M1", "Thbs1", "Itga4", "M1", 3, 1, 524.5768889999999, 1573.730667, 0.5894633891046084, 0.6088102578564109, 3, 1, 49.39115366666666, 148.173461, 0.4028566641859401, 0.4047308569616938, 25909.45773458094, 233185.1196112284, 0.2374692545944213, 0.2464042973892949),
    @("M1", "Thbs1", "Itga4", "M2", 3, 1, 524.5768889999999, 1573.730667, 0.5894633891046084, 0.6088102578564109, 3, 1, 70.865733, 212.597199, 0.5780130789036135, 0.5807021443531359, 37174.52575384464, 334570.7317846017, 0.3407175484373134, 0.3535374222414034),
    @("M1", "Thbs1", "Itga4", "Neutro", 3, 1, 524.5768889999999, 1573.730667, 0.5894633891046084, 0.6088102578564109, 3, 1, 0.6137786666666667, 1.841336, 0.005006257352694774, 0.005029547749002215, 321.974103494568, 2897.766931451112, 0.002951005425849326, 0.003062040261971169),
    @("M1", "Thbs1", "Itga4", "sCs", 3, 1, 524.5768889999999, 1573.730667, 0.5894633891046084, 0.6088102578564109, 2, 1, 1.3859455, 2.771891, 0.01130440047304115, 0.007571327633592945, 727.0349787135494, 4362.209872281297, 0.006663530214634576, 0.00460950192892309),
    @("M2", "Thbs1", "Itga4", "ECs", 3, 1, 181.4813383333333, 544.4440149999999, 0.2039293133121744, 0.210622508737405, 2, 1, 0.317263, 0.634526, 0.002587741009497455, 0.001733186563985812, 57.57731384364833, 345.46388306189, 0.000527716247096569, 0.0003650481022166546),
    @("M2", "Thbs1", "Itga4", "FAPs", 3, 1, 181.4813383333333, 544.4440149999999, 0.2039293133121744, 0.210622508737405, 1, 0.3333333333333333, 0.02842633333333333, 0.085279, 0.0002318580752130288, 0.0002329367385893503, 5.158849017242776, 46.42964115518499, 0.00004728265806407544, 0.00004906172025879804),
    @("M2", "Thbs1", "Itga4", "M1", 3, 1, 181.4813383333333, 544.4440149999999, 0.2039293133121744, 0.210622508737405, 3, 1, 49.39115366666666, 148.173461, 0.4028566641859401, 0.4047308569616938, 8963.572669253988, 80672.1540232859, 0.082154282890672, 0.08524542845671174),
    @("M2", "Thbs1", "Itga4", "M2", 3, 1, 181.4813383333333, 544.4440149999999, 0.2039293133121744, 0.210622508737405, 3, 1, 70.865733, 212.597199, 0.5780130789036135, 0.5807021443531359, 12860.80806681266, 115747.272601314, 0.1178738102662696, 0.1223089424728482),
    @("M2", "Thbs1", "Itga4", "Neutro", 3, 1, 181.4813383333333, 544.4440149999999, 0.2039293133121744, 0.210622508737405, 3, 1, 0.6137786666666667, 1.841336, 0.005006257352694774, 0.005029547749002215, 111.3893738671155, 1002.50436480404, 0.001020922624199069, 0.001059335964709414),
    @("M2", "Thbs1", "Itga4", "sCs", 3, 1, 181.4813383333333, 544.4440149999999, 0.2039293133121744, 0.210622508737405, 2, 1, 1.3859455, 2.771891, 0.01130440047304115, 0.007571327633592945, 251.5232441970608, 1509.139465182365, 0.002305298625873102, 0.001594692020660186),
    @("Neutro", "Thbs1", "Itga4", "ECs", 3, 1, 15.79677433333333, 47.390323, 0.01775072507139627, 0.01833332435500452, 2, 1, 0.317263, 0.634526, 0.002587741009497455, 0.001733186563985812, 5.011732015316333, 30.070392091898, 0.00004593427921556677, 0.00003177507144528769),
    @("Neutro", "Thbs1", "Itga4", "FAPs", 3, 1, 15.79677433333333, 47.390323, 0.01775072507139627, 0.01833332435500452, 1, 0.3333333333333333, 0.02842633333333333, 0.085279, 0.0002318580752130288, 0.0002329367385893503, 0.4490443727907777, 4.041399355117, 0.000004115648948689592, 0.000004270504782755456),
    @("Neutro", "Thbs1", "Itga4", "M1", 3, 1, 15.79677433333333, 47.390323, 0.01775072507139627, 0.01833332435500452, 3, 1, 49.39115366666666, 148.173461, 0.4028566641859401, 0.4047308569616938, 780.2209085353225, 7021.988176817902, 0.007150997889144434, 0.00742006207715767),
    @("Neutro", "Thbs1", "Itga4", "M2", 3, 1, 15.79677433333333, 47.390323, 0.01775072507139627, 0.01833332435500452, 3, 1, 70.865733, 212.597199, 0.5780130789036135, 0.5807021443531359, 1119.449992167253, 10075.04992950528, 0.01026015125128932, 0.0106462007660727),
    @("Neutro", "Thbs1", "Itga4", "Neutro", 3, 1, 15.79677433333333, 47.390323, 0.01775072507139627, 0.01833332435500452, 3, 1, 0.6137786666666667, 1.841336, 0.005006257352694774, 0.005029547749002215, 9.695723087947556, 87.261507791528, 0.00008886469790434103, 0.00009220833024144046),
    @("Neutro", "Thbs1", "Itga4", "sCs", 3, 1, 15.79677433333333, 47.390323, 0.01775072507139627, 0.01833332435500452, 2, 1, 1.3859455, 2.771891, 0.01130440047304115, 0.007571327633592945, 21.89346830179883, 131.360809810793, 0.0002006613048939154, 0.0001388076053046682),
    @("sCs", "Thbs1", "Itga4", "ECs", 2, 1, 53.1948375, 106.389675, 0.05977466764766092, 0.0411576941520005, 2, 1, 0.317263, 0.634526, 0.002587741009497455, 0.001733186563985812, 16.8767537297625, 67.50701491905001, 0.0001546813588009329, 0.00007133396250888472),
    @("sCs", "Thbs1", "Itga4", "FAPs", 2, 1, 53.1948375, 106.389675, 0.05977466764766092, 0.0411576941520005, 1, 0.3333333333333333, 0.02842633333333333, 0.085279, 0.0002318580752130288, 0.0002329367385893503, 1.5121341823875, 9.072805094324998, 0.00001385923938728516, 0.000009587139043624972),
    @("sCs", "Thbs1", "Itga4", "M1", 2, 1, 53.1948375, 106.389675, 0.05977466764766092, 0.0411576941520005, 3, 1, 49.39115366666666, 148.173461, 0.4028566641859401, 0.4047308569616938, 2627.354393235862, 15764.12635941517, 0.02408062321135991, 0.01665778882470646),
    @("sCs", "Thbs1", "Itga4", "M2", 2, 1, 53.1948375, 106.389675, 0.05977466764766092, 0.0411576941520005, 3, 1, 70.865733, 212.597199, 0.5780130789036135, 0.5807021443531359, 3769.691151253387, 22618.14690752033, 0.03455053968746471, 0.02390036125069721),
    @("sCs", "Thbs1", "Itga4", "Neutro", 2, 1, 53.1948375, 106.389675, 0.05977466764766092, 0.0411576941520005, 3, 1, 0.6137786666666667, 1.841336, 0.005006257352694774, 0.005029547749002215, 32.6498564343, 195.8991386058, 0.0002992473694159889, 0.0002070045879763158),
    @("sCs", "Thbs1", "Itga4", "sCs", 2, 1, 53.1948375, 106.389675, 0.05977466764766092, 0.0411576941520005, 2, 1, 1.3859455, 2.771891, 0.01130440047304115, 0.007571327633592945, 73.72514565635625, 294.900582625425, 0.0006757167812320958, 0.0003116183870680082),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
